{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" paragraphs, along\n// with the single blank paragraph that immediately precedes them, leaving\n// the rest of the document (including the blank paragraph + page break that\n// follow) untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two target paragraphs by their exact text.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"Powered by Jekyll and Github pages\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx !== -1 && copyrightIdx !== -1) {\n  // The blank paragraph right before \"Ver no Jupiter...\" is also removed.\n  const blankIdx = jupiterIdx - 1;\n  if (blankIdx >= 0 && items[blankIdx].text === \"\") {\n    items[blankIdx].delete();\n  }\n  items[jupiterIdx].delete();\n  items[copyrightIdx].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" paragraphs, along\n# with the single blank paragraph that immediately precedes them, leaving the\n# rest of the document (including the blank paragraph + page break that\n# follow) untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the start of the \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n# paragraph, then step back one paragraph so the deletion also swallows the\n# blank paragraph right before it.\n$startRange = $d.Content\n$found1 = $startRange.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\nif ($found1) {\n    $startRange.Collapse(1) | Out-Null  # wdCollapseStart\n    $startRange.Move(4, -1) | Out-Null  # wdParagraph, back one paragraph (the blank line)\n    $delStart = $startRange.Start\n\n    # Locate the copyright paragraph (\"Powered by Jekyll\" is a stable, plain\n    # ASCII substring of it) and expand to its full paragraph extent so the\n    # deletion also removes its trailing paragraph mark.\n    $endRange = $d.Content\n    $found2 = $endRange.Find.Execute(\"Powered by Jekyll\")\n    if ($found2) {\n        $endRange.Expand(4) | Out-Null  # wdParagraph\n        $delEnd = $endRange.End\n\n        $d.Range($delStart, $delEnd).Delete() | Out-Null\n    }\n}\n"}
